$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'61.035.19"
$ws.Range("E2").Value = "  -1.76%  "

# Row 3
$ws.Range("D3").Value = "'2.972.65"
$ws.Range("E3").Value = "  -0.39%  "

# Row 4
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.08%  "

# Row 5
$ws.Range("D5").Value = "'596.72"
$ws.Range("E5").Value = "  +3.06%  "

# Row 6
$ws.Range("D6").Value = "'142.62"
$ws.Range("E6").Value = "  -1.87%  "

# Row 7
$ws.Range("E7").Value = "  +0.03%  "

# Row 8
$ws.Range("D8").Value = "'0.515"
$ws.Range("E8").Value = "  -0.90%  "

# Row 9
$ws.Range("D9").Value = "'2.969.87"
$ws.Range("E9").Value = "  -0.40%  "

# Row 10
$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").Value = "'0.145"
$ws.Range("E10").Value = "  -1.48%  "

# Row 11
$ws.Range("B11").Value = "Toncoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D11").Value = "'5.98"
$ws.Range("E11").Value = "  +6.34%  "

# Row 12
$ws.Range("D12").Value = "'0.453"
$ws.Range("E12").Value = "  +3.12%  "

# Row 13
$ws.Range("D13").Value = "'0.0000225"
$ws.Range("E13").Value = "  -0.19%  "

# Row 14
$ws.Range("D14").Value = "'34.03"
$ws.Range("E14").Value = "  -1.03%  "

# Row 15
$ws.Range("D15").Value = "'0.124"
$ws.Range("E15").Value = "  +2.36%  "

# Row 16
$ws.Range("D16").Value = "'3.471.22"
$ws.Range("E16").Value = "  -0.12%  "

# Row 17
$ws.Range("D17").Value = "'6.86"
$ws.Range("E17").Value = "  -1.94%  "

# Row 18
$ws.Range("D18").Value = "'61.045.33"
$ws.Range("E18").Value = "  -1.71%  "

# Row 19
$ws.Range("D19").Value = "'2.974.21"
$ws.Range("E19").Value = "  -0.56%  "

# Row 20
$ws.Range("D20").Value = "'446.16"
$ws.Range("E20").Value = "  -1.67%  "

# Row 21
$ws.Range("D21").Value = "'14.10"
$ws.Range("E21").Value = "  +2.24%  "

# Row 22
$ws.Range("D22").Value = "'0.678"
$ws.Range("E22").Value = "  +0.71%  "

# Row 23
$ws.Range("D23").Value = "'7.26"
$ws.Range("E23").Value = "  +0.16%  "

# Row 24
$ws.Range("D24").Value = "'81.84"
$ws.Range("E24").Value = "  +2.67%  "

# Row 25
$ws.Range("D25").Value = "'2.16"
$ws.Range("E25").Value = "  -4.77%  "

# Row 26
$ws.Range("D26").Value = "'10.35"
$ws.Range("E26").Value = "  +3.87%  "

# Row 27
$ws.Range("D27").Value = "'11.83"
$ws.Range("E27").Value = "  -2.75%  "

# Row 28
$ws.Range("E28").Value = "  -0.05%  "

# Row 29
$ws.Range("D29").Value = "'2.67"
$ws.Range("E29").Value = "  +3.12%  "

# Row 30
$ws.Range("E30").Value = "  -0.09%  "

# Row 31
$ws.Range("D31").Value = "'7.08"
$ws.Range("E31").Value = "  -1.32%  "

# Row 32
$ws.Range("D32").Value = "'2.03"
$ws.Range("E32").Value = "  -2.08%  "

# Row 33
$ws.Range("D33").Value = "'27.02"
$ws.Range("E33").Value = "  +1.25%  "

# Row 35
$ws.Range("D35").Value = "'0.0₃0808"
$ws.Range("E35").Value = "  +3.90%  "

# Row 36
$ws.Range("E36").Value = "  -0.56%  "

# Row 37
$ws.Range("D37").Value = "'5.72"
$ws.Range("E37").Value = "  +0.43%  "

# Row 38
$ws.Range("D38").Value = "'50.21"
$ws.Range("E38").Value = "  +0.48%  "

# Row 39
$ws.Range("E39").Value = "  -2.29%  "

# Row 40
$ws.Range("D40").Value = "'8.99"
$ws.Range("E40").Value = "  +0.53%  "

# Row 41
$ws.Range("E41").Value = "  +9.54%  "

# Row 42
$ws.Range("E42").Value = "  -2.28%  "

# Row 43
$ws.Range("D43").Value = "'391.14"
$ws.Range("E43").Value = "  -3.54%  "

# Row 44
$ws.Range("D44").Value = "'39.03"
$ws.Range("E44").Value = "  +1.89%  "

# Row 45
$ws.Range("D45").Value = "'0.0347"
$ws.Range("E45").Value = "  -0.35%  "

# Row 46
$ws.Range("D46").Value = "'0.265"
$ws.Range("E46").Value = "  -3.68%  "

# Row 47
$ws.Range("D47").Value = "'2.679.42"
$ws.Range("E47").Value = "  -2.75%  "

# Row 48
$ws.Range("D48").Value = "'130.46"
$ws.Range("E48").Value = "  +2.52%  "

# Row 50
$ws.Range("E50").Value = "  -0.55%  "

# Row 51
$ws.Range("E51").Value = "  -0.24%  "
